$d = $word.ActiveDocument

$replacements = @(
    @{old = "40×29=1160"; new = "51×84=4284"},
    @{old = "25×22=550"; new = "76×59=4484"},
    @{old = "74×59=4366"; new = "40×62=2480"},
    @{old = "72×23=1656"; new = "35×65=2275"},
    @{old = "11×20=220"; new = "14×65=910"},
    @{old = "59×12=708"; new = "61×23=1403"},
    @{old = "23×72=1656"; new = "34×79=2686"},
    @{old = "18×60=1080"; new = "40×94=3760"},
    @{old = "57×60=3420"; new = "73×76=5548"},
    @{old = "22×80=1760"; new = "27×47=1269"},
    @{old = "53×73=3869"; new = "76×99=7524"},
    @{old = "41×82=3362"; new = "22×50=1100"},
    @{old = "24×29=696"; new = "89×74=6586"},
    @{old = "14×90=1260"; new = "73×77=5621"},
    @{old = "60×83=4980"; new = "36×24=864"},
    @{old = "75×99=7425"; new = "49×92=4508"},
    @{old = "24×96=2304"; new = "81×37=2997"},
    @{old = "67×61=4087"; new = "24×14=336"},
    @{old = "25×80=2000"; new = "46×74=3404"},
    @{old = "69×53=3657"; new = "48×44=2112"},
    @{old = "15×13=195"; new = "81×42=3402"},
    @{old = "94×24=2256"; new = "24×72=1728"},
    @{old = "32×48=1536"; new = "18×14=252"},
    @{old = "51×60=3060"; new = "15×35=525"},
    @{old = "78×21=1638"; new = "32×91=2912"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $r.new, 2)
}
